$d = $word.ActiveDocument

# Update the date paragraph at the top of the document
$d.Paragraphs.Item(1).Range.Text = "2023-11-26 Sunday"

# Update each answer cell in the 20x5 table, row-major order
$t = $d.Tables.Item(1)
$answers = @(
  "99-12=87",
  "50-17=33",
  "48-32=16",
  "91-48=43",
  "83-75=8",
  "40-1=39",
  "35+25=60",
  "41-8=33",
  "29-18=11",
  "68+31=99",
  "93-62=31",
  "73+8=81",
  "94-51=43",
  "97-5=92",
  "17-6=11",
  "38+32=70",
  "72+5=77",
  "42-41=1",
  "66-49=17",
  "50+17=67",
  "76-25=51",
  "90-12=78",
  "85+0=85",
  "2+10=12",
  "31+41=72",
  "22+37=59",
  "14+69=83",
  "80-14=66",
  "94-47=47",
  "67+25=92",
  "63-49=14",
  "17+11=28",
  "28+23=51",
  "14+8=22",
  "0+78=78",
  "63+18=81",
  "39+39=78",
  "81-24=57",
  "11+80=91",
  "39-11=28",
  "34+56=90",
  "75-10=65",
  "17+76=93",
  "28+19=47",
  "91-59=32",
  "56-30=26",
  "94-32=62",
  "35-13=22",
  "21-19=2",
  "48+4=52",
  "41+30=71",
  "8+43=51",
  "10-10=0",
  "30+67=97",
  "67-27=40",
  "70-14=56",
  "42+33=75",
  "28+14=42",
  "59-12=47",
  "36+54=90",
  "13+69=82",
  "9+8=17",
  "93-60=33",
  "91-80=11",
  "4+56=60",
  "83-71=12",
  "53+41=94",
  "40+5=45",
  "43-12=31",
  "38+49=87",
  "55+12=67",
  "17+72=89",
  "94-56=38",
  "23+37=60",
  "79-40=39",
  "10+50=60",
  "51-2=49",
  "14+7=21",
  "3+79=82",
  "56-44=12",
  "98-64=34",
  "84+6=90",
  "59-39=20",
  "18-5=13",
  "19+31=50",
  "43+18=61",
  "57-23=34",
  "92-19=73",
  "40+30=70",
  "66+1=67",
  "80-72=8",
  "84-37=47",
  "7+17=24",
  "32-1=31",
  "73+17=90",
  "78-78=0",
  "11+86=97",
  "80-36=44",
  "31+0=31",
  "36+35=71"
)

$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
  for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $answers[$idx]
    $idx = $idx + 1
  }
}

Write-Host "Done. Updated $idx cells."
